$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Drop the old "_GoBack" bookmark that previously sat in the empty
#    paragraph right after the equipment table.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Turn the "Cedula: XXXXXXXXXX" placeholder into "Cedula: [cedula]",
#    split across three runs ("[", "cedula", "]") that drop w:eastAsia /
#    w:szCs, with a freshly minted "_GoBack" bookmark sitting right before
#    the closing "]" run.
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Cedula:*") {
        $target = $p
        break
    }
}

$pRange = $target.Range
# Exclude the trailing paragraph mark so the replacement stays inside the
# same paragraph (keeps w:pPr / paraId / rsids untouched).
$contentRange = $d.Range($pRange.Start, $pRange.End - 1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="42C9255A" w14:textId="68D6BD96" w:rsidR="004B053B" w:rsidRDefault="004B053B" w:rsidP="004B053B">
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="1440"/>
              </w:tabs>
              <w:spacing w:line="360" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve">Cedula: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:sz w:val="24"/>
              </w:rPr>
              <w:t>[</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:sz w:val="24"/>
              </w:rPr>
              <w:t>cedula</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:sz w:val="24"/>
              </w:rPr>
              <w:t>]</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

[void]$contentRange.InsertXML($xml)
